$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scraper bug fix: stray commas inside multi-party "Razon social" names were
#     ambiguous with the CSV-style comma separator used by the scraper; replace with periods. ---
$ws.Range("E37").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E91").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E131").Value = 'RICCOTTI. MARIANA EDITH'
$ws.Range("E183").Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range("F145").Value = 'MERCANZINI. GASTON ARIEL'

# --- Scraper bug fix: "Importe" amounts were scraped in Spanish locale format
#     ("1.234,56" - "." thousands separator, "," decimal separator). Re-save them as
#     plain "1234.56" text (no thousands separator, "." decimal separator). The column
#     is marked as Text so the literal string (with trailing zeroes) survives. ---
$importe = @{
    2 = '7000.00'
    3 = '36300.00'
    4 = '75000.00'
    5 = '36000.00'
    6 = '255000.00'
    7 = '4085.00'
    8 = '50.00'
    9 = '46845.00'
    10 = '1500.00'
    11 = '960.00'
    12 = '4840.00'
    13 = '3180.00'
    14 = '800.00'
    15 = '24680.00'
    16 = '512793.20'
    17 = '106100.36'
    18 = '4000.00'
    19 = '10909.84'
    20 = '7880.00'
    21 = '9036.00'
    22 = '4434.39'
    23 = '1575.00'
    24 = '16872.24'
    25 = '6000.00'
    26 = '3825.00'
    27 = '4895.00'
    28 = '23340.00'
    29 = '50.45'
    30 = '68435.74'
    31 = '68.00'
    32 = '38405.00'
    33 = '48227.37'
    34 = '1147.68'
    35 = '2365.00'
    36 = '211.48'
    37 = '960.00'
    38 = '2280.08'
    39 = '165210.00'
    40 = '566.00'
    41 = '9.87'
    42 = '130.10'
    43 = '73.15'
    44 = '96765.00'
    45 = '3715.53'
    46 = '2700.00'
    47 = '39227.49'
    48 = '2400.00'
    49 = '4098.60'
    50 = '62.50'
    51 = '5561.09'
    52 = '10141.00'
    53 = '586.96'
    54 = '4352.00'
    55 = '91.00'
    56 = '234.30'
    57 = '2150.00'
    58 = '16.03'
    59 = '138.50'
    60 = '15068.70'
    61 = '133575.00'
    62 = '41203.11'
    63 = '4388.00'
    64 = '69.32'
    65 = '1275.00'
    66 = '450.00'
    67 = '18640.90'
    68 = '1367.06'
    69 = '2889.07'
    70 = '3034.14'
    71 = '1200.00'
    72 = '2639.00'
    73 = '8000.00'
    74 = '1950.00'
    75 = '520.00'
    76 = '400.00'
    77 = '3728.50'
    78 = '840.00'
    79 = '241323.69'
    80 = '42000.00'
    81 = '3817.25'
    82 = '3997.00'
    83 = '1836.00'
    84 = '28100.00'
    85 = '2080.00'
    86 = '880.00'
    87 = '2320.50'
    88 = '15260.00'
    89 = '1437.00'
    90 = '15894.00'
    91 = '33720.00'
    92 = '2267.19'
    93 = '129.00'
    94 = '73980.00'
    95 = '20.00'
    96 = '3880.00'
    97 = '1300.00'
    98 = '3000.00'
    99 = '391178.31'
    100 = '53261.49'
    101 = '2875.00'
    102 = '59.81'
    103 = '116.66'
    104 = '32230.79'
    105 = '3597.00'
    106 = '23500.00'
    107 = '4615.20'
    108 = '45.28'
    109 = '9153.00'
    110 = '5091.78'
    111 = '1489.00'
    112 = '3335.93'
    113 = '1880.00'
    114 = '20270.81'
    115 = '3502.60'
    116 = '8767.60'
    117 = '66.00'
    118 = '100.00'
    119 = '30.00'
    120 = '2850.00'
    121 = '23797.79'
    122 = '29.00'
    123 = '6795.00'
    124 = '82.97'
    125 = '8354.75'
    126 = '9500.00'
    127 = '12570.00'
    128 = '1800.00'
    129 = '378.00'
    130 = '2148.00'
    131 = '2000.00'
    132 = '27000.00'
    133 = '850.00'
    134 = '2300.00'
    135 = '188500.00'
    136 = '1758.00'
    137 = '49215.06'
    138 = '18978.00'
    139 = '9600.00'
    140 = '780.00'
    141 = '1200.00'
    142 = '4700.00'
    143 = '4000.00'
    144 = '2785.00'
    145 = '9000.00'
    146 = '7000.00'
    147 = '1155.60'
    148 = '240.00'
    149 = '81.50'
    150 = '283.00'
    151 = '12225.00'
    152 = '91.30'
    153 = '707886.32'
    154 = '9303.46'
    155 = '2300.00'
    156 = '7000.00'
    157 = '4000.00'
    158 = '30026.15'
    159 = '6292.00'
    160 = '1800.00'
    161 = '1657.50'
    162 = '3312.00'
    163 = '4000.00'
    164 = '800.00'
    165 = '1000.00'
    166 = '32700.00'
    167 = '12000.00'
    168 = '2400.00'
    169 = '3000.00'
    170 = '6042.40'
    171 = '3000.00'
    172 = '600.00'
    173 = '1400.00'
    174 = '8400.00'
    175 = '3900.00'
    176 = '760.00'
    177 = '3900.00'
    178 = '1500.00'
    179 = '11590.00'
    180 = '143.32'
    181 = '1400.00'
    182 = '1800.00'
    183 = '6070.00'
    184 = '99.70'
    185 = '22860.00'
    186 = '2686.00'
    187 = '432.00'
    188 = '194.66'
    189 = '115.00'
    190 = '13688.00'
    191 = '8494.79'
    192 = '1567.76'
    193 = '11028.00'
    194 = '23.26'
    195 = '6560.00'
    196 = '1802.00'
    197 = '8160.00'
    198 = '4714.27'
    199 = '456.40'
    200 = '473.68'
    201 = '1270.80'
    202 = '150.00'
    203 = '948.00'
    204 = '8800.00'
    205 = '1540.00'
    206 = '500.00'
    207 = '1200.00'
    208 = '5100.00'
    209 = '372401.12'
    210 = '2781.60'
    211 = '2161.76'
    212 = '1694600.70'
    213 = '2300.00'
    214 = '9440.00'
    215 = '1344850.00'
    216 = '354000.00'
    217 = '1017500.00'
    218 = '1458630.00'
    219 = '1071931.00'
    220 = '32500.00'
    221 = '204000.00'
    222 = '927950.00'
    223 = '2380832.00'
    224 = '1318300.00'
    225 = '130000.00'
    226 = '120800.00'
    227 = '958500.00'
    228 = '463900.00'
    229 = '550000.00'
    230 = '84000.00'
    231 = '3988.00'
    232 = '13980.00'
    233 = '8300.00'
    234 = '46200.00'
    235 = '1000.00'
    236 = '58800.00'
    237 = '37.00'
    238 = '5900.00'
    239 = '449100.00'
    240 = '1823.99'
    241 = '400.00'
    242 = '4500.00'
    243 = '5000.00'
}

$ws.Range("H2:H243").NumberFormat = "@"
foreach ($row in $importe.Keys) {
    $ws.Range("H$row").Value = $importe[$row]
}
